$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cart_Page")

# --- New row: qty-update locator added under "update buttons" ---
$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "update buttons"
$ws.Range("B15").Value = "class name"
$ws.Range("C15").Value = "cart-qty-update"

# --- New row: delete-item locator added under "delete buttons" ---
$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = "delete buttons"
$ws.Range("B17").Value = "class name"
$ws.Range("C17").Value = "itemDelete-CartPage"

# --- New row: item subtotal value locator ---
$ws.Rows("21:21").Insert()
$ws.Range("A21").Value = "item subtotals"
$ws.Range("B21").Value = "class name"
$ws.Range("C21").Value = "cart-item-value.total-value"

# --- New row: cart total (grand total) locator ---
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "cart total"
$ws.Range("B20").Value = "class name"
$ws.Range("C20").Value = "cart-total-value.cart-total-grandTotal"

# --- New row: cart item title block locator, inserted at the top of this block ---
$ws.Rows("5:5").Insert()
$ws.Range("A5").Value = "drop down products"
$ws.Range("B5").Value = "class name"
$ws.Range("C5").Value = "cart-item-block.cart-item-title"

# --- Remove the now-superfluous extra trailing blank row ---
$ws.Rows("25:25").Delete()

# --- Update the selection shown in the sheet view ---
$ws.Activate()
$ws.Range("A4:B5").Select()
